$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.37788999080658
$ws.Range("B1").Value = 6.415459632873535
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 3.689490079879761
$ws.Range("E1").Value = 2.013293504714966
